$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are stored as text so formatted values
# (including trailing zeros, e.g. "0.06200") are preserved exactly,
# matching the source workbook which stores prices as text.

# --- Simple price updates (column D) ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "276.51"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.16"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.267"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06200"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "6.572"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1662"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08308"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03515"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03168"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09152"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.759"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001639"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006278"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006220"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001068"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.717"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.313"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3292"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002735"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04758"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01132"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006269"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7226"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001399"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01239"

# --- Row reorder: coins at rows 41-43 shifted down one position ---
# (KickToken/BKEXToken/CEJI -> CEJI/KickToken/BKEXToken, with updated prices)

$ws.Range("B41").Value = "CEJI"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005297"
$ws.Range("E41").Value = "40CEJICEJI"

$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007071"
$ws.Range("E42").Value = "41KickTokenKICK"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1122"
$ws.Range("E43").Value = "42BKEXTokenBKK"
